$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '29.034.38'
Set-TextValue $ws.Range('D3') '1.829.13'
$ws.Range('E3').Value = '  -0.36%  '
Set-TextValue $ws.Range('D4') '0.9990'
$ws.Range('E4').Value = '  +0.00%  '
Set-TextValue $ws.Range('D5') '241.22'
$ws.Range('E5').Value = '  -0.37%  '
Set-TextValue $ws.Range('D6') '0.6272'
$ws.Range('E6').Value = '  -5.16%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('B8').Value = 'OKB'
$ws.Range('C8').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range('D8') '45.12'
$ws.Range('E8').Value = '  +1.18%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Range('D9') '0.07590'
$ws.Range('E9').Value = '  +2.15%  '
$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue $ws.Range('D10') '0.2912'
$ws.Range('E10').Value = '  -0.88%  '
$ws.Range('B11').Value = 'Solana'
$ws.Range('C11').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue $ws.Range('D11') '22.77'
$ws.Range('E11').Value = '  -0.67%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range('D12') '0.07641'
$ws.Range('E12').Value = '  -1.38%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D13') '1.830.79'
$ws.Range('E13').Value = '  +0.16%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D14') '4.954'
$ws.Range('E14').Value = '  -0.92%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range('D15') '0.6646'
$ws.Range('E15').Value = '  -0.38%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range('D16') '82.33'
$ws.Range('E16').Value = '  -0.93%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range('D17') '0.000009424'
$ws.Range('E17').Value = '  +9.72%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range('D18') '5.980'
$ws.Range('E18').Value = '  -2.34%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range('D19') '28.862.34'
$ws.Range('E19').Value = '  -1.08%  '
Set-TextValue $ws.Range('D20') '224.82'
$ws.Range('E20').Value = '  -0.93%  '
$ws.Range('E21').Value = '  -1.25%  '
$ws.Range('E22').Value = '  -0.18%  '
Set-TextValue $ws.Range('D23') '7.224'
$ws.Range('E23').Value = '  +1.68%  '
Set-TextValue $ws.Range('D25') '160.21'
$ws.Range('E25').Value = '  +0.16%  '
Set-TextValue $ws.Range('D26') '8.410'
$ws.Range('E26').Value = '  -2.60%  '
Set-TextValue $ws.Range('D27') '0.1364'
$ws.Range('E27').Value = '  -2.74%  '
Set-TextValue $ws.Range('D28') '17.81'
$ws.Range('E28').Value = '  -0.89%  '
Set-TextValue $ws.Range('D29') '1.495'
$ws.Range('E29').Value = '  -1.41%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D30') '4.045'
$ws.Range('E30').Value = '  -1.69%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D31') '4.023'
$ws.Range('E31').Value = '  -0.51%  '
$ws.Range('E32').Value = '  +1.18%  '
Set-TextValue $ws.Range('D33') '0.05197'
$ws.Range('E33').Value = '  -1.69%  '
Set-TextValue $ws.Range('D34') '1.841'
$ws.Range('E34').Value = '  -1.46%  '
$ws.Range('E35').Value = '  +0.41%  '
Set-TextValue $ws.Range('D36') '0.7292'
$ws.Range('E36').Value = '  -1.08%  '
Set-TextValue $ws.Range('D37') '2.611'
$ws.Range('E37').Value = '  -1.85%  '
Set-TextValue $ws.Range('D38') '1.271.90'
$ws.Range('E38').Value = '  -2.03%  '
$ws.Range('E39').Value = '  +0.82%  '
Set-TextValue $ws.Range('D40') '0.01784'
$ws.Range('E40').Value = '  -0.58%  '
Set-TextValue $ws.Range('D41') '6.489'
$ws.Range('E41').Value = '  +7.22%  '
Set-TextValue $ws.Range('D42') '0.8888'
$ws.Range('E42').Value = '  -3.35%  '
$ws.Range('E43').Value = '  +0.06%  '
Set-TextValue $ws.Range('D44') '101.38'
$ws.Range('E44').Value = '  -0.82%  '
Set-TextValue $ws.Range('D45') '1.973.38'
$ws.Range('E45').Value = '  -0.28%  '
Set-TextValue $ws.Range('D46') '0.5103'
$ws.Range('E46').Value = '  -0.69%  '
$ws.Range('E47').Value = '  +0.24%  '
Set-TextValue $ws.Range('D49') '0.3977'
$ws.Range('E49').Value = '  -0.86%  '
Set-TextValue $ws.Range('D50') '0.07300'
$ws.Range('E50').Value = '  -13.24%  '
Set-TextValue $ws.Range('D51') '8.794'
$ws.Range('E51').Value = '  +0.83%  '
